# Applies the "Validate Feedback 2 @Validate" commit to
# ProductionGradeSpringBootApi.xlsx
#
# Summary of changes:
#  - "2.Validation" sheet gets a new feedback row (row 5, mirroring the
#    existing row 4 "Validate Feedback 1" pair) plus three new rows appended
#    after the existing row 18 entry (rows 19-21) describing the "Validate
#    Feedback 2" follow up / solution.
#  - Column D on that sheet is widened to fit the new, longer text.
#  - Selection/active-cell bookkeeping is updated on the "Design" sheet and
#    the "2.Validation" sheet to match where the author left the cursor.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Design" sheet: author's selection moved from C15 to F4
# ---------------------------------------------------------------------
$wsDesign = $wb.Worksheets.Item("Design")
$wsDesign.Activate() | Out-Null
$wsDesign.Range("F4").Select() | Out-Null

# ---------------------------------------------------------------------
# "2.Validation" sheet: new feedback content
# ---------------------------------------------------------------------
$wsValidation = $wb.Worksheets.Item("2.Validation")
$wsValidation.Activate() | Out-Null

# New row 5 mirrors row 4 ("Validate Feedback 1 @Valid" / "1. Use @Valid ...")
# but for the second feedback round.
$wsValidation.Range("D5").Value = "Validate Feedback 2 @Validate : "
$wsValidation.Range("E5").Value = "2. Use @Validate and provide custom message "

# New rows appended below the existing row 18 entry.
$wsValidation.Range("C19").Value = 2
$wsValidation.Range("D19").Value = "Use @Validate and provide custom message and POST check input should not be null and call custom application validation"
$wsValidation.Range("D20").Value = "Solution: HandlerInterceptor.prehandle"
$wsValidation.Range("C21").Value = 3

# Widen column D (no longer just best-fit) to accommodate the new text.
$wsValidation.Columns.Item(4).ColumnWidth = 28.8

# Author left the cursor on D10 after editing.
$wsValidation.Range("D10").Select() | Out-Null
